$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.27
$ws.Range("E2").Value = 1.14
$ws.Range("B3").Value = 1.57
$ws.Range("D3").Value = 1.37
$ws.Range("D4").Value = 1.33
$ws.Range("G4").Value = 1
$ws.Range("B5").Value = 1.6
$ws.Range("E6").Value = 1.34
$ws.Range("D7").Value = 1.71
